$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The comments table (ID / Comments / Time) grows by one more submitted
# entry: a new row right after the last used row, continuing the existing
# ID sequence (35 -> 36).
$lastRow = $ws.Range("A1048576").End(-4162).Row   # xlUp
$newRow = $lastRow + 1

$ws.Range("A" + $newRow).Value = 36
$ws.Range("B" + $newRow).Value = "dsfdsfsd"

# Column C stores dates as plain text (e.g. "04-11-2023"), so force text
# formatting before assigning the value - otherwise Excel would interpret
# the date-shaped string and convert it to a date serial number.
$ws.Range("C" + $newRow).NumberFormat = "@"
$ws.Range("C" + $newRow).Value = "04-11-2023"
